$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 6365.933
$ws.Cells.Item(116, 9).Value = 5104.9
$ws.Cells.Item(116, 11).Value = 5104.9
$ws.Cells.Item(116, 13).Value = -1662.9
$ws.Cells.Item(132, 8).Value = 1262.69
$ws.Cells.Item(132, 9).Value = 1277.551
$ws.Cells.Item(132, 11).Value = 3832.653
$ws.Cells.Item(132, 13).Value = -1302.653
$ws.Cells.Item(137, 8).Value = 8266.808000000001
$ws.Cells.Item(137, 9).Value = 4173.5
$ws.Cells.Item(137, 10).Value = 12884.897
$ws.Cells.Item(137, 11).Value = 12520.5
$ws.Cells.Item(137, 12).Value = 38654.69100000001
$ws.Cells.Item(137, 13).Value = -9970.5
$ws.Cells.Item(137, 14).Value = -43754.69100000001
$ws.Cells.Item(138, 8).Value = 13025.245
$ws.Cells.Item(138, 9).Value = 23499
$ws.Cells.Item(138, 10).Value = 12614.51
$ws.Cells.Item(138, 11).Value = 70497
$ws.Cells.Item(138, 12).Value = 37843.53
$ws.Cells.Item(138, 13).Value = -65357
$ws.Cells.Item(138, 14).Value = -48123.53

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 16959.617
$ws.Cells.Item(61, 9).Value = 4002.95
$ws.Cells.Item(61, 11).Value = 4002.95
$ws.Cells.Item(61, 13).Value = -3790.95
$ws.Cells.Item(74, 8).Value = 12335.878
$ws.Cells.Item(74, 9).Value = 2098
$ws.Cells.Item(74, 11).Value = 2098
$ws.Cells.Item(74, 13).Value = -1224
$ws.Cells.Item(77, 8).Value = 12335.878
$ws.Cells.Item(77, 9).Value = 2098
$ws.Cells.Item(77, 11).Value = 10490
$ws.Cells.Item(77, 13).Value = -6122
$ws.Cells.Item(122, 8).Value = 4953.1562
$ws.Cells.Item(122, 9).Value = 2843.8823
$ws.Cells.Item(122, 10).Value = 7343.6665
$ws.Cells.Item(122, 11).Value = 8531.6469
$ws.Cells.Item(122, 12).Value = 22030.9995
$ws.Cells.Item(122, 13).Value = -6081.6469
$ws.Cells.Item(122, 14).Value = -26930.9995
$ws.Cells.Item(125, 8).Value = 59933.332
$ws.Cells.Item(125, 10).Value = 59933.332
$ws.Cells.Item(125, 12).Value = 59933.332
$ws.Cells.Item(125, 14).Value = -69773.33199999999
$ws.Cells.Item(132, 8).Value = 10771.556
$ws.Cells.Item(132, 9).Value = 4646.483
$ws.Cells.Item(132, 10).Value = 36146.855
$ws.Cells.Item(132, 11).Value = 13939.449
$ws.Cells.Item(132, 12).Value = 108440.565
$ws.Cells.Item(132, 13).Value = -11409.449
$ws.Cells.Item(132, 14).Value = -113500.565
$ws.Cells.Item(136, 8).Value = 16959.617
$ws.Cells.Item(136, 9).Value = 4002.95
$ws.Cells.Item(136, 11).Value = 12008.85
$ws.Cells.Item(136, 13).Value = -9458.849999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(64, 8).Value = 787.5833
$ws.Cells.Item(64, 9).Value = 612
$ws.Cells.Item(64, 10).Value = 875.375
$ws.Cells.Item(64, 11).Value = 612
$ws.Cells.Item(64, 12).Value = 875.375
$ws.Cells.Item(64, 13).Value = -387
$ws.Cells.Item(64, 14).Value = -1325.375
$ws.Cells.Item(67, 8).Value = 787.5833
$ws.Cells.Item(67, 9).Value = 612
$ws.Cells.Item(67, 10).Value = 875.375
$ws.Cells.Item(67, 11).Value = 612
$ws.Cells.Item(67, 12).Value = 875.375
$ws.Cells.Item(67, 13).Value = 168
$ws.Cells.Item(67, 14).Value = -2435.375
$ws.Cells.Item(86, 8).Value = 8875
$ws.Cells.Item(86, 9).Value = 7833.3335
$ws.Cells.Item(86, 11).Value = 7833.3335
$ws.Cells.Item(86, 13).Value = -6710.3335
$ws.Cells.Item(89, 8).Value = 8875
$ws.Cells.Item(89, 9).Value = 7833.3335
$ws.Cells.Item(89, 11).Value = 39166.6675
$ws.Cells.Item(89, 13).Value = -33550.6675
$ws.Cells.Item(94, 8).Value = 2911.28
$ws.Cells.Item(94, 9).Value = 3073.0527
$ws.Cells.Item(94, 10).Value = 2399
$ws.Cells.Item(94, 11).Value = 3073.0527
$ws.Cells.Item(94, 12).Value = 2399
$ws.Cells.Item(94, 13).Value = -2622.0527
$ws.Cells.Item(94, 14).Value = -3301
$ws.Cells.Item(99, 8).Value = 3240.05
$ws.Cells.Item(99, 9).Value = 3450.0588
$ws.Cells.Item(99, 11).Value = 3450.0588
$ws.Cells.Item(99, 13).Value = -1952.0588
$ws.Cells.Item(107, 8).Value = 1860.84
$ws.Cells.Item(107, 9).Value = 1840.0952
$ws.Cells.Item(107, 10).Value = 1969.75
$ws.Cells.Item(107, 11).Value = 1840.0952
$ws.Cells.Item(107, 12).Value = 1969.75
$ws.Cells.Item(107, 13).Value = 79.90480000000002
$ws.Cells.Item(107, 14).Value = -5809.75
$ws.Cells.Item(134, 8).Value = 10494.333
$ws.Cells.Item(134, 9).Value = 5099.926
$ws.Cells.Item(134, 11).Value = 15299.778
$ws.Cells.Item(134, 13).Value = -12764.778

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 2721.3
$ws.Cells.Item(22, 9).Value = 1959
$ws.Cells.Item(22, 11).Value = 1959
$ws.Cells.Item(22, 13).Value = -1609
$ws.Cells.Item(31, 8).Value = 17411.277
$ws.Cells.Item(31, 9).Value = 8485.200000000001
$ws.Cells.Item(31, 11).Value = 8485.200000000001
$ws.Cells.Item(31, 13).Value = -8190.200000000001
$ws.Cells.Item(34, 8).Value = 17411.277
$ws.Cells.Item(34, 9).Value = 8485.200000000001
$ws.Cells.Item(34, 11).Value = 8485.200000000001
$ws.Cells.Item(34, 13).Value = -8283.200000000001
$ws.Cells.Item(107, 8).Value = 307923.4
$ws.Cells.Item(107, 9).Value = 479042.12
$ws.Cells.Item(107, 11).Value = 479042.12
$ws.Cells.Item(107, 13).Value = -477122.12
$ws.Cells.Item(122, 8).Value = 5535.6943
$ws.Cells.Item(122, 9).Value = 2476.8635
$ws.Cells.Item(122, 11).Value = 7430.5905
$ws.Cells.Item(122, 13).Value = -4980.5905
$ws.Cells.Item(127, 8).Value = 76768
$ws.Cells.Item(127, 10).Value = 76768
$ws.Cells.Item(127, 12).Value = 76768
$ws.Cells.Item(127, 14).Value = -86688
$ws.Cells.Item(134, 8).Value = 5324.657
$ws.Cells.Item(134, 9).Value = 1502.2222
$ws.Cells.Item(134, 11).Value = 4506.6666
$ws.Cells.Item(134, 13).Value = -1971.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 5405.4
$ws.Cells.Item(114, 10).Value = 7999.5
$ws.Cells.Item(114, 12).Value = 23998.5
$ws.Cells.Item(114, 14).Value = -30506.5
$ws.Cells.Item(122, 8).Value = 15377110
$ws.Cells.Item(122, 9).Value = 62289564
$ws.Cells.Item(122, 10).Value = 2582804.2
$ws.Cells.Item(122, 11).Value = 560606076
$ws.Cells.Item(122, 12).Value = 23245237.8
$ws.Cells.Item(122, 13).Value = -560603626
$ws.Cells.Item(122, 14).Value = -23250137.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 4639.439
$ws.Cells.Item(132, 9).Value = 5275.1377
$ws.Cells.Item(132, 11).Value = 15825.4131
$ws.Cells.Item(132, 13).Value = -13295.4131

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 3210.55
$ws.Cells.Item(16, 9).Value = 3532.7856
$ws.Cells.Item(16, 10).Value = 2458.6667
$ws.Cells.Item(16, 11).Value = 3532.7856
$ws.Cells.Item(16, 12).Value = 2458.6667
$ws.Cells.Item(16, 13).Value = -3362.7856
$ws.Cells.Item(16, 14).Value = -2798.6667
$ws.Cells.Item(68, 8).Value = 8521.888999999999
$ws.Cells.Item(68, 10).Value = 10113.857
$ws.Cells.Item(68, 12).Value = 10113.857
$ws.Cells.Item(68, 14).Value = -11611.857
$ws.Cells.Item(71, 8).Value = 8521.888999999999
$ws.Cells.Item(71, 10).Value = 10113.857
$ws.Cells.Item(71, 12).Value = 50569.285
$ws.Cells.Item(71, 14).Value = -58057.285
$ws.Cells.Item(93, 8).Value = 13232.6
$ws.Cells.Item(93, 9).Value = 6610.4443
$ws.Cells.Item(93, 11).Value = 6610.4443
$ws.Cells.Item(93, 13).Value = -5362.4443
$ws.Cells.Item(101, 8).Value = 24330
$ws.Cells.Item(101, 10).Value = 24330
$ws.Cells.Item(101, 12).Value = 24330
$ws.Cells.Item(101, 14).Value = -30820
$ws.Cells.Item(136, 8).Value = 13761.018
$ws.Cells.Item(136, 9).Value = 10595.263
$ws.Cells.Item(136, 11).Value = 31785.789
$ws.Cells.Item(136, 13).Value = -29235.789

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 4303.625
$ws.Cells.Item(2, 9).Value = 7775
$ws.Cells.Item(2, 10).Value = 832.25
$ws.Cells.Item(2, 11).Value = 7775
$ws.Cells.Item(2, 12).Value = 832.25
$ws.Cells.Item(2, 13).Value = -7663
$ws.Cells.Item(2, 14).Value = -1056.25
$ws.Cells.Item(123, 8).Value = 35000
$ws.Cells.Item(123, 10).Value = 35000
$ws.Cells.Item(123, 12).Value = 35000
$ws.Cells.Item(123, 14).Value = -44800
$ws.Cells.Item(132, 8).Value = 8680.852999999999
$ws.Cells.Item(132, 9).Value = 3637.3044
$ws.Cells.Item(132, 10).Value = 19226.455
$ws.Cells.Item(132, 11).Value = 10911.9132
$ws.Cells.Item(132, 12).Value = 57679.36500000001
$ws.Cells.Item(132, 13).Value = -8381.913199999999
$ws.Cells.Item(132, 14).Value = -62739.36500000001
$ws.Cells.Item(136, 8).Value = 10098.593
$ws.Cells.Item(136, 9).Value = 2126.3333
$ws.Cells.Item(136, 10).Value = 38001.5
$ws.Cells.Item(136, 11).Value = 6378.999899999999
$ws.Cells.Item(136, 12).Value = 114004.5
$ws.Cells.Item(136, 13).Value = -3828.999899999999
$ws.Cells.Item(136, 14).Value = -119104.5
